# Update the "想去人数" (want-to-go count) values in column F
# for the 展览 (sheet1) and 全部类型 (sheet4) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 173
$wsExhibit.Range("F5").Value = 233
$wsExhibit.Range("F8").Value = 2257
$wsExhibit.Range("F10").Value = 5568

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 173
$wsAll.Range("F6").Value = 233
$wsAll.Range("F11").Value = 2257
$wsAll.Range("F13").Value = 5568
